$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Current Row/Column" property values in column B
$ws.Range("B2").Value = 15
$ws.Range("B3").Value = 17
$ws.Range("B5").Value = 6
